$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the old rows 2-6 into a single Python-tuple-like string in A2
$ws.Range("A2").Value = "('False Prophet', ['{2}{W}{W}', 'Creature " + [char]0x2014 + " Human Cleric', 'When False Prophet dies, exile all creatures.', '2/2'])"

# Remove the now-obsolete rows 3-6 entirely so the used range shrinks to A1:A2
$ws.Range("A3:A6").EntireRow.Delete()
